# Raul's Log.xlsx - "still in local mode" update
# - Rows 424-430: fix date 42622 -> 42621 and add the "Arrive 10 minutes
#   early..." special-instructions note in column F.
# - Rows 431-442: fix date 42622 -> 42621 (F column already populated).
# - Rows 447-458: append new log entries for 8/29 (serial 42625).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$arriveNote = "Arrive 10 minutes early. Ensure that the instructor does not require further assistance before you leave."

# ---- Rows 424-430: date correction + new special instructions note ----
foreach ($r in 424..430) {
    $ws.Range("B$r").Value = 42621
    $ws.Range("F$r").Value = $arriveNote
    $ws.Range("F$r").WrapText = $true
    $ws.Range("F$r").HorizontalAlignment = -4108
    $ws.Rows.Item($r).RowHeight = 30
}

# ---- Rows 431-442: date correction only ----
foreach ($r in 431..442) {
    $ws.Range("B$r").Value = 42621
}

# ---- New rows 447-458 ----

# Row 447
$ws.Range("A447").Value = "Pickup Skype Kit"
$ws.Range("B447").Value = 42625
$ws.Range("C447").Value = "1700"
$ws.Range("D447").Value = "HNE"
$ws.Range("E447").Value = "402"

# Row 448
$ws.Range("A448").Value = "AV Shutdown"
$ws.Range("B448").Value = 42625
$ws.Range("C448").Value = "1700"
$ws.Range("D448").Value = "HNE"
$ws.Range("E448").Value = "402"

# Row 449
$ws.Range("A449").Value = "Demo"
$ws.Range("B449").Value = 42625
$ws.Range("C449").Value = "1630"
$ws.Range("D449").Value = "HNE"
$ws.Range("E449").Value = "B11"

# Row 450
$ws.Range("A450").Value = "AV Shutdown"
$ws.Range("B450").Value = 42625
$ws.Range("C450").Value = "1830"
$ws.Range("D450").Value = "HNE"
$ws.Range("E450").Value = "B11"

# Row 451
$ws.Range("A451").Value = "Demo"
$ws.Range("B451").Value = 42625
$ws.Range("C451").Value = "1900"
$ws.Range("D451").Value = "DB"
$ws.Range("E451").Value = "0004"

# Row 452
$ws.Range("A452").Value = "Demo"
$ws.Range("B452").Value = 42625
$ws.Range("C452").Value = "1900"
$ws.Range("D452").Value = "DB"
$ws.Range("E452").Value = "0016"

# Row 453
$ws.Range("A453").Value = "Demo"
$ws.Range("B453").Value = 42625
$ws.Range("C453").Value = "1900"
$ws.Range("D453").Value = "HNE"
$ws.Range("E453").Value = "038"

# Row 454
$ws.Range("A454").Value = "Demo"
$ws.Range("B454").Value = 42625
$ws.Range("C454").Value = "1900"
$ws.Range("D454").Value = "HNE"
$ws.Range("E454").Value = "103"

# Row 455
$ws.Range("A455").Value = "Demo"
$ws.Range("B455").Value = 42625
$ws.Range("C455").Value = "1900"
$ws.Range("D455").Value = "HNE"
$ws.Range("E455").Value = "401"

# Row 456 (wrapped F note, row height 45)
$ws.Range("A456").Value = "Demo"
$ws.Range("B456").Value = 42625
$ws.Range("C456").Value = "1630"
$ws.Range("D456").Value = "SSB"
$ws.Range("E456").Value = "W141"
$ws.Range("F456").Value = "Podium mic (there ) and 3 neck mics (2 built in - located in drawer;  third one to be plugged into mixer In rear booth- there) . Test all mics and demo to client"
$ws.Range("F456").WrapText = $true
$ws.Range("F456").HorizontalAlignment = -4108
$ws.Rows.Item(456).RowHeight = 45

# Row 457 (wrapped F note, default row height)
$ws.Range("A457").Value = "Operator"
$ws.Range("B457").Value = 42625
$ws.Range("C457").Value = "1700"
$ws.Range("D457").Value = "SSB"
$ws.Range("E457").Value = "W141"
$ws.Range("F457").Value = "Operate event between 17:00 - 18:00"
$ws.Range("F457").WrapText = $true
$ws.Range("F457").HorizontalAlignment = -4108

# Row 458 (wrapped F note, row height 30)
$ws.Range("A458").Value = "Pickup Mic"
$ws.Range("B458").Value = 42625
$ws.Range("C458").Value = "2000"
$ws.Range("D458").Value = "SSB"
$ws.Range("E458").Value = "W141"
$ws.Range("F458").Value = "Return 2 neck mics to drawer and third to rear booth; leave podium mic in place"
$ws.Range("F458").WrapText = $true
$ws.Range("F458").HorizontalAlignment = -4108
$ws.Rows.Item(458).RowHeight = 30

# ---- View state: select A458 (matches the author scrolling to the new
# last row after appending entries) ----
$ws.Range("A458").Select()
